# Weekly refresh of Fruta / Vega Modelo de Temuco - Coco price data.
# The diff re-shuffles the Fecha/Volumen/Precio values across rows 2-49
# (a like-for-like permutation of the same underlying weekly records),
# so we just overwrite each changed cell with its new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44475
$ws.Range("M2").Value = 20

# Row 3
$ws.Range("D3").Value = 44251

# Row 4
$ws.Range("D4").Value = 44414
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 25000
$ws.Range("S4").Value = 1250

# Row 5
$ws.Range("D5").Value = 44434
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 24000
$ws.Range("P5").Value = 24000
$ws.Range("S5").Value = 1200

# Row 6
$ws.Range("D6").Value = 44419
$ws.Range("M6").Value = 40

# Row 7
$ws.Range("D7").Value = 44452
$ws.Range("M7").Value = 25
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 1250

# Row 8
$ws.Range("D8").Value = 44391
$ws.Range("M8").Value = 10

# Row 9
$ws.Range("D9").Value = 44433
$ws.Range("M9").Value = 10

# Row 10
$ws.Range("D10").Value = 44466
$ws.Range("M10").Value = 70
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("S10").Value = 1200

# Row 11
$ws.Range("D11").Value = 44425
$ws.Range("M11").Value = 15
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("S11").Value = 1200

# Row 12
$ws.Range("D12").Value = 44454
$ws.Range("M12").Value = 25
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 25000
$ws.Range("S12").Value = 1250

# Row 13
$ws.Range("D13").Value = 44426
$ws.Range("M13").Value = 15

# Row 14
$ws.Range("D14").Value = 44421
$ws.Range("M14").Value = 20

# Row 15
$ws.Range("D15").Value = 44467

# Row 16
$ws.Range("D16").Value = 44235
$ws.Range("M16").Value = 15
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 25000
$ws.Range("S16").Value = 1250

# Row 17
$ws.Range("D17").Value = 44334
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 25000
$ws.Range("O17").Value = 25000
$ws.Range("P17").Value = 25000
$ws.Range("S17").Value = 1250

# Row 18
$ws.Range("D18").Value = 44356
$ws.Range("M18").Value = 15
$ws.Range("N18").Value = 24000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 24000
$ws.Range("S18").Value = 1200

# Row 19
$ws.Range("D19").Value = 44175
$ws.Range("M19").Value = 25
$ws.Range("N19").Value = 23000
$ws.Range("O19").Value = 23000
$ws.Range("P19").Value = 23000
$ws.Range("S19").Value = 1150

# Row 20
$ws.Range("D20").Value = 44222

# Row 21
$ws.Range("D21").Value = 44238
$ws.Range("M21").Value = 30

# Row 22
$ws.Range("D22").Value = 44468

# Row 23
$ws.Range("D23").Value = 44442
$ws.Range("M23").Value = 25
$ws.Range("N23").Value = 23000
$ws.Range("O23").Value = 23000
$ws.Range("P23").Value = 23000
$ws.Range("S23").Value = 1150

# Row 24
$ws.Range("D24").Value = 44214
$ws.Range("M24").Value = 15
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 25000
$ws.Range("P24").Value = 25000
$ws.Range("S24").Value = 1250

# Row 25
$ws.Range("D25").Value = 44389
$ws.Range("M25").Value = 20
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 24000
$ws.Range("S25").Value = 1200

# Row 26
$ws.Range("D26").Value = 44349
$ws.Range("M26").Value = 30
$ws.Range("N26").Value = 24000
$ws.Range("O26").Value = 24000
$ws.Range("P26").Value = 24000
$ws.Range("S26").Value = 1200

# Row 27
$ws.Range("D27").Value = 44412
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 25000
$ws.Range("S27").Value = 1250

# Row 28
$ws.Range("D28").Value = 44398
$ws.Range("M28").Value = 15
$ws.Range("N28").Value = 25000
$ws.Range("O28").Value = 25000
$ws.Range("P28").Value = 25000
$ws.Range("S28").Value = 1250

# Row 29
$ws.Range("D29").Value = 44420
$ws.Range("M29").Value = 35

# Row 30
$ws.Range("D30").Value = 44249
$ws.Range("M30").Value = 15

# Row 31
$ws.Range("D31").Value = 44232
$ws.Range("N31").Value = 25000
$ws.Range("O31").Value = 25000
$ws.Range("P31").Value = 25000
$ws.Range("S31").Value = 1250

# Row 32
$ws.Range("D32").Value = 44431

# Row 33
$ws.Range("D33").Value = 44400
$ws.Range("M33").Value = 5
$ws.Range("N33").Value = 24000
$ws.Range("O33").Value = 24000
$ws.Range("P33").Value = 24000
$ws.Range("S33").Value = 1200

# Row 34
$ws.Range("D34").Value = 44363
$ws.Range("M34").Value = 30

# Row 35
$ws.Range("D35").Value = 44390
$ws.Range("M35").Value = 10

# Row 36
$ws.Range("D36").Value = 44461
$ws.Range("M36").Value = 30

# Row 37
$ws.Range("D37").Value = 44221
$ws.Range("M37").Value = 30
$ws.Range("N37").Value = 25000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 25000
$ws.Range("S37").Value = 1250

# Row 38
$ws.Range("D38").Value = 44432
$ws.Range("M38").Value = 30

# Row 39
$ws.Range("D39").Value = 44428
$ws.Range("N39").Value = 24000
$ws.Range("O39").Value = 24000
$ws.Range("P39").Value = 24000
$ws.Range("S39").Value = 1200

# Row 40
$ws.Range("D40").Value = 44462
$ws.Range("M40").Value = 10

# Row 41
$ws.Range("D41").Value = 44435
$ws.Range("M41").Value = 100

# Row 42
$ws.Range("D42").Value = 44231

# Row 43
$ws.Range("D43").Value = 44489
$ws.Range("M43").Value = 40

# Row 44
$ws.Range("D44").Value = 44396
$ws.Range("M44").Value = 12
$ws.Range("N44").Value = 24000
$ws.Range("O44").Value = 24000
$ws.Range("P44").Value = 24000
$ws.Range("S44").Value = 1200

# Row 45
$ws.Range("D45").Value = 44418
$ws.Range("M45").Value = 20

# Row 46
$ws.Range("D46").Value = 44469
$ws.Range("M46").Value = 40

# Row 47
$ws.Range("D47").Value = 44392
$ws.Range("M47").Value = 10
$ws.Range("N47").Value = 24000
$ws.Range("O47").Value = 24000
$ws.Range("P47").Value = 24000
$ws.Range("S47").Value = 1200

# Row 48
$ws.Range("D48").Value = 44474
$ws.Range("M48").Value = 20

# Row 49
$ws.Range("D49").Value = 44424
$ws.Range("M49").Value = 25

